# Append the 2025-04-29 price row (row 59) to every price-history sheet in
# the workbook, duplicating the last known (2025-04-28 / row 58) price for
# each series — this mirrors how the daily Argent-price feed appends one
# new row per day.

$wb = $excel.ActiveWorkbook

# Sheet name -> new price value for column B, row 59.
$newPrices = [ordered]@{
    "N-Dense"                   = "38"
    "N-Type"                    = "37.33"
    "N-type Wafer"              = "1.13"
    "Cell Topcon 183mm"         = "0.275"
    "Module Topcon 183mm"       = "0.09"
    "Silver Rear_side"          = "5,342"
    "Silver Busbar front-side"  = "7,997"
    "Silver finger front-side"  = "8,047"
    "USD_CNY"                   = "7.3078"
}

$newDate = "2025-04-29"
$newRow = 59

foreach ($sheetName in $newPrices.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    $dateCell = $ws.Cells.Item($newRow, 1)
    $priceCell = $ws.Cells.Item($newRow, 2)

    # Force plain text storage so "2025-04-29" isn't auto-coerced into a
    # date serial, matching the rest of the column which stores the date
    # (and price) as literal text.
    $dateCell.NumberFormat = "@"
    $priceCell.NumberFormat = "@"

    $dateCell.Value = $newDate
    $priceCell.Value = $newPrices[$sheetName]
}
